# Weekly CompStat update: new crime data collected (19th Precinct)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: issue number and reporting week dates ---
$ws.Range("A8").Value = 'Volume 32   Number  41'
$ws.Range("C9").Value = 'Report Covering the Week  10/6/2025  Through  10/12/2025'

# --- Crime complaint statistics table ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 9
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -40
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -18.181818181818
$ws.Range("N15").Value = -64

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 15
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -6.25
$ws.Range("I16").Value = 135
$ws.Range("J16").Value = 168
$ws.Range("K16").Value = -19.642857142857
$ws.Range("L16").Value = -16.666666666666
$ws.Range("M16").Value = 22.727272727272
$ws.Range("N16").Value = -87.394957983193

# Row 17
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = -4.761904761904
$ws.Range("I17").Value = 160
$ws.Range("J17").Value = 149
$ws.Range("K17").Value = 7.38255033557
$ws.Range("L17").Value = 2.564102564102
$ws.Range("M17").Value = 107.792207792208
$ws.Range("N17").Value = -31.914893617021

# Row 18
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 26
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 136.363636363636
$ws.Range("I18").Value = 229
$ws.Range("J18").Value = 175
$ws.Range("K18").Value = 30.857142857142
$ws.Range("L18").Value = 9.56937799043
$ws.Range("M18").Value = 21.164021164021
$ws.Range("N18").Value = -89.808633733867

# Row 19
$ws.Range("C19").Value = 29
$ws.Range("D19").Value = 44
$ws.Range("E19").Value = -34.090909090909
$ws.Range("F19").Value = 115
$ws.Range("G19").Value = 157
$ws.Range("H19").Value = -26.751592356687
$ws.Range("I19").Value = 1291
$ws.Range("J19").Value = 1286
$ws.Range("K19").Value = 0.388802488335
$ws.Range("L19").Value = -4.157386785449
$ws.Range("M19").Value = 32.139201637666
$ws.Range("N19").Value = -55.344171566931

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 0
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 6
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 44
$ws.Range("J20").Value = 72
$ws.Range("K20").Value = -38.888888888888
$ws.Range("L20").Value = -66.412213740458
$ws.Range("M20").Value = -37.142857142857
$ws.Range("N20").Value = -98.357596117954

# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 61
$ws.Range("E21").Value = -29.508196721311
$ws.Range("F21").Value = 184
$ws.Range("G21").Value = 216
$ws.Range("H21").Value = -14.814814814814
$ws.Range("I21").Value = 1868
$ws.Range("J21").Value = 1867
$ws.Range("K21").Value = 0.053561863952
$ws.Range("L21").Value = -7.38720872583
$ws.Range("M21").Value = 30.174216027874
$ws.Range("N21").Value = -79.59139080083

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 24
$ws.Range("J22").Value = 36
$ws.Range("K22").Value = -33.333333333333

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 26
$ws.Range("K23").Value = 4
$ws.Range("L23").Value = 18.181818181818
$ws.Range("M23").Value = 18.181818181818

# Row 24
$ws.Range("C24").Value = 52
$ws.Range("D24").Value = 70
$ws.Range("E24").Value = -25.714285714285
$ws.Range("F24").Value = 217
$ws.Range("G24").Value = 296
$ws.Range("H24").Value = -26.689189189189
$ws.Range("I24").Value = 2294
$ws.Range("J24").Value = 2556
$ws.Range("K24").Value = -10.250391236306
$ws.Range("L24").Value = -8.203281312525
$ws.Range("M24").Value = 71.321882001493

# Row 25
$ws.Range("C25").Value = 36
$ws.Range("D25").Value = 71
$ws.Range("E25").Value = -49.295774647887
$ws.Range("F25").Value = 154
$ws.Range("G25").Value = 269
$ws.Range("H25").Value = -42.750929368029
$ws.Range("I25").Value = 1891
$ws.Range("J25").Value = 2242
$ws.Range("K25").Value = -15.655664585191
$ws.Range("L25").Value = -13.494967978042

# Row 26
$ws.Range("C26").Value = 12
$ws.Range("E26").Value = 71.428571428571
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = 25
$ws.Range("I26").Value = 286
$ws.Range("J26").Value = 284
$ws.Range("K26").Value = 0.704225352112
$ws.Range("L26").Value = 7.924528301886
$ws.Range("M26").Value = 4.379562043795

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 12
$ws.Range("J27").Value = 20
$ws.Range("K27").Value = -40
$ws.Range("L27").Value = -40

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = 42.857142857142
$ws.Range("I28").Value = 76
$ws.Range("J28").Value = 93
$ws.Range("K28").Value = -18.279569892473
$ws.Range("L28").Value = 7.042253521126

# Row 31
$ws.Range("D31").Value = 1
$ws.Range("D31").NumberFormat = "#,##0"
$ws.Range("E31").Value = -100
$ws.Range("E31").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J31").Value = 29
$ws.Range("K31").Value = -48.275862068965
$ws.Range("L31").Value = -6.25

# --- Cells reverting to the "0" placeholder text (no incidents reported) ---
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
